# Update Wnt1-Ror2.xlsx with new TPM-derived values.
# The data now only contains the "FAPs" sending cluster (rows for the
# "ECs" sending cluster are removed), and the remaining rows/columns are
# refreshed with new expression/specificity figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows for the second ("ECs") sending cluster block
# (originally rows 5-7), shrinking the sheet from 7 rows to 4 rows.
$ws.Rows("5:7").Delete() | Out-Null

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt1"
$ws.Range("C2").Value = "Ror2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3532066666666667
$ws.Range("H2").Value = 1.05962
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04009133333333333
$ws.Range("N2").Value = 0.120274
$ws.Range("O2").Value = 0.01033409631432067
$ws.Range("P2").Value = 0.01033409631432067
$ws.Range("Q2").Value = 0.01416052620888889
$ws.Range("R2").Value = 0.12744473588
$ws.Range("S2").Value = 0.01033409631432067
$ws.Range("T2").Value = 0.01033409631432067

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Ror2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3532066666666667
$ws.Range("H3").Value = 1.05962
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.244570333333333
$ws.Range("N3").Value = 6.733711
$ws.Range("O3").Value = 0.578569084147867
$ws.Range("P3").Value = 0.578569084147867
$ws.Range("Q3").Value = 0.7927972055355555
$ws.Range("R3").Value = 7.135174849819999
$ws.Range("S3").Value = 0.578569084147867
$ws.Range("T3").Value = 0.578569084147867

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt1"
$ws.Range("C4").Value = "Ror2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3532066666666667
$ws.Range("H4").Value = 1.05962
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.594858333333333
$ws.Range("N4").Value = 4.784575
$ws.Range("O4").Value = 0.4110968195378122
$ws.Range("P4").Value = 0.4110968195378122
$ws.Range("Q4").Value = 0.5633145957222223
$ws.Range("R4").Value = 5.0698313615
$ws.Range("S4").Value = 0.4110968195378122
$ws.Range("T4").Value = 0.4110968195378122
